# Add season-record columns (Wins, Losses, Ties) to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header style from an existing header cell (A1) onto the new headers
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill season record for each data row (2-53)
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 72   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 90   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
